$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new label in A13 (appends a new shared string "Commit test")
$ws.Range("A13").Value = "Commit test"

# Move the active selection to A14 (next empty row below the table)
[void]$ws.Range("A14").Select()
